$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2102.5908
$ws.Range("J28").Value = 6601
$ws.Range("L28").Value = 6601
$ws.Range("N28").Value = -7571
$ws.Range("H33").Value = 223.625
$ws.Range("I33").Value = 241.85715
$ws.Range("K33").Value = 241.85715
$ws.Range("M33").Value = -12.85714999999999
$ws.Range("H43").Value = 2508.6
$ws.Range("I43").Value = 2907.1667
$ws.Range("J43").Value = 2337.7856
$ws.Range("K43").Value = 2907.1667
$ws.Range("L43").Value = 2337.7856
$ws.Range("M43").Value = -2838.1667
$ws.Range("N43").Value = -2475.7856
$ws.Range("H55").Value = 643.63635
$ws.Range("J55").Value = 626.4286
$ws.Range("L55").Value = 626.4286
$ws.Range("N55").Value = -1054.4286
$ws.Range("H62").Value = 14173.6
$ws.Range("I62").Value = 17216
$ws.Range("K62").Value = 17216
$ws.Range("M62").Value = -16592
$ws.Range("H65").Value = 14173.6
$ws.Range("I65").Value = 17216
$ws.Range("K65").Value = 86080
$ws.Range("M65").Value = -82960
$ws.Range("H106").Value = 2827.2942
$ws.Range("I106").Value = 2768.8572
$ws.Range("K106").Value = 2768.8572
$ws.Range("M106").Value = -2137.8572
$ws.Range("H107").Value = 784
$ws.Range("I107").Value = 730.375
$ws.Range("K107").Value = 730.375
$ws.Range("M107").Value = 1189.625
$ws.Range("H116").Value = 2325000.2
$ws.Range("I116").Value = 2325000.2
$ws.Range("K116").Value = 2325000.2
$ws.Range("M116").Value = -2321558.2
$ws.Range("H131").Value = 1112526.6
$ws.Range("I131").Value = 1112526.6
$ws.Range("K131").Value = 3337579.8
$ws.Range("M131").Value = -3332539.8
$ws.Range("H135").Value = 459.82352
$ws.Range("I135").Value = 329.85715
$ws.Range("K135").Value = 2968.71435
$ws.Range("M135").Value = -433.7143499999997
$ws.Range("H138").Value = 261601.97
$ws.Range("I138").Value = 4061.5625
$ws.Range("J138").Value = 388391.1
$ws.Range("K138").Value = 12184.6875
$ws.Range("L138").Value = 1165173.3
$ws.Range("M138").Value = -7044.6875
$ws.Range("N138").Value = -1175453.3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3748.0225
$ws.Range("I32").Value = 2687.6506
$ws.Range("K32").Value = 2687.6506
$ws.Range("M32").Value = -2400.6506
$ws.Range("H45").Value = 33112.93
$ws.Range("I45").Value = 43765.3
$ws.Range("K45").Value = 43765.3
$ws.Range("M45").Value = -43388.3
$ws.Range("H61").Value = 5053.231
$ws.Range("I61").Value = 1426.3572
$ws.Range("K61").Value = 1426.3572
$ws.Range("M61").Value = -1214.3572
$ws.Range("H74").Value = 170830.06
$ws.Range("I74").Value = 266304.1
$ws.Range("K74").Value = 266304.1
$ws.Range("M74").Value = -265430.1
$ws.Range("H77").Value = 170830.06
$ws.Range("I77").Value = 266304.1
$ws.Range("K77").Value = 1331520.5
$ws.Range("M77").Value = -1327152.5
$ws.Range("H110").Value = 4142.136
$ws.Range("I110").Value = 3575.125
$ws.Range("K110").Value = 3575.125
$ws.Range("M110").Value = -1530.125
$ws.Range("H132").Value = 2559.5278
$ws.Range("I132").Value = 2116.4443
$ws.Range("J132").Value = 3888.7778
$ws.Range("K132").Value = 6349.3329
$ws.Range("L132").Value = 11666.3334
$ws.Range("M132").Value = -3819.3329
$ws.Range("N132").Value = -16726.3334
$ws.Range("H136").Value = 5053.231
$ws.Range("I136").Value = 1426.3572
$ws.Range("K136").Value = 4279.071599999999
$ws.Range("M136").Value = -1729.071599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2695.634
$ws.Range("I134").Value = 2426.5667
$ws.Range("J134").Value = 3429.4546
$ws.Range("K134").Value = 7279.7001
$ws.Range("L134").Value = 10288.3638
$ws.Range("M134").Value = -4744.7001
$ws.Range("N134").Value = -15358.3638
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1703.3572
$ws.Range("I16").Value = 1804.2727
$ws.Range("K16").Value = 1804.2727
$ws.Range("M16").Value = -1517.2727
$ws.Range("H31").Value = 4093.6326
$ws.Range("I31").Value = 3412.2173
$ws.Range("J31").Value = 4696.423
$ws.Range("K31").Value = 3412.2173
$ws.Range("L31").Value = 4696.423
$ws.Range("M31").Value = -3117.2173
$ws.Range("N31").Value = -5286.423
$ws.Range("H34").Value = 4093.6326
$ws.Range("I34").Value = 3412.2173
$ws.Range("J34").Value = 4696.423
$ws.Range("K34").Value = 3412.2173
$ws.Range("L34").Value = 4696.423
$ws.Range("M34").Value = -3210.2173
$ws.Range("N34").Value = -5100.423
$ws.Range("H58").Value = 4534.9473
$ws.Range("I58").Value = 4457.6665
$ws.Range("K58").Value = 4457.6665
$ws.Range("M58").Value = -4254.6665
$ws.Range("H59").Value = 92017
$ws.Range("J59").Value = 92017
$ws.Range("L59").Value = 92017
$ws.Range("N59").Value = -94307
$ws.Range("H107").Value = 517.1515000000001
$ws.Range("I107").Value = 523.5599999999999
$ws.Range("K107").Value = 523.5599999999999
$ws.Range("M107").Value = 1396.44
$ws.Range("H113").Value = 1703.3572
$ws.Range("I113").Value = 1804.2727
$ws.Range("K113").Value = 1804.2727
$ws.Range("M113").Value = 365.7273
$ws.Range("H122").Value = 3727.5715
$ws.Range("I122").Value = 3204.5
$ws.Range("J122").Value = 4425
$ws.Range("K122").Value = 9613.5
$ws.Range("L122").Value = 13275
$ws.Range("M122").Value = -7163.5
$ws.Range("N122").Value = -18175
$ws.Range("H132").Value = 2436.2456
$ws.Range("I132").Value = 1903
$ws.Range("J132").Value = 4435.9165
$ws.Range("K132").Value = 5709
$ws.Range("L132").Value = 13307.7495
$ws.Range("M132").Value = -3179
$ws.Range("N132").Value = -18367.7495
$ws.Range("H134").Value = 2184.8438
$ws.Range("I134").Value = 2169.6
$ws.Range("K134").Value = 6508.799999999999
$ws.Range("M134").Value = -3973.799999999999
$ws.Range("H135").Value = 119270.55
$ws.Range("J135").Value = 119270.55
$ws.Range("L135").Value = 119270.55
$ws.Range("N135").Value = -129410.55
$ws.Range("H136").Value = 4534.9473
$ws.Range("I136").Value = 4457.6665
$ws.Range("K136").Value = 13372.9995
$ws.Range("M136").Value = -10822.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5061.909
$ws.Range("I131").Value = 12279.923
$ws.Range("J131").Value = 2035
$ws.Range("K131").Value = 36839.769
$ws.Range("L131").Value = 6105
$ws.Range("M131").Value = -31799.769
$ws.Range("N131").Value = -16185
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2728.9092
$ws.Range("I122").Value = 2701.8
$ws.Range("K122").Value = 8105.400000000001
$ws.Range("M122").Value = -5655.400000000001
$ws.Range("H126").Value = 4425.278
$ws.Range("I126").Value = 1977.3636
$ws.Range("J126").Value = 8272
$ws.Range("K126").Value = 5932.0908
$ws.Range("L126").Value = 24816
$ws.Range("M126").Value = -3462.0908
$ws.Range("N126").Value = -29756
$ws.Range("H132").Value = 2320.7903
$ws.Range("I132").Value = 1961.6666
$ws.Range("J132").Value = 3271.4119
$ws.Range("K132").Value = 5884.9998
$ws.Range("L132").Value = 9814.235700000001
$ws.Range("M132").Value = -3354.9998
$ws.Range("N132").Value = -14874.2357
$ws.Range("H135").Value = 68248.086
$ws.Range("J135").Value = 68248.086
$ws.Range("L135").Value = 68248.086
$ws.Range("N135").Value = -78388.086
$ws.Range("H136").Value = 12522.3125
$ws.Range("J136").Value = 12522.3125
$ws.Range("L136").Value = 37566.9375
$ws.Range("N136").Value = -42666.9375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1923.25
$ws.Range("I61").Value = 1843.4546
$ws.Range("K61").Value = 1843.4546
$ws.Range("M61").Value = -1641.4546
$ws.Range("H100").Value = 1043625.1
$ws.Range("I100").Value = 1355483
$ws.Range("K100").Value = 1355483
$ws.Range("M100").Value = -1354942
$ws.Range("H113").Value = 1923.25
$ws.Range("I113").Value = 1843.4546
$ws.Range("K113").Value = 1843.4546
$ws.Range("M113").Value = 326.5454
$ws.Range("H122").Value = 2684.3333
$ws.Range("I122").Value = 2004
$ws.Range("K122").Value = 6012
$ws.Range("M122").Value = -3562
$ws.Range("H136").Value = 4921
$ws.Range("I136").Value = 4554.875
$ws.Range("K136").Value = 13664.625
$ws.Range("M136").Value = -11114.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 12332.833
$ws.Range("I74").Value = 15984
$ws.Range("J74").Value = 10507.25
$ws.Range("K74").Value = 15984
$ws.Range("L74").Value = 10507.25
$ws.Range("M74").Value = -15048
$ws.Range("N74").Value = -12379.25
$ws.Range("H77").Value = 12332.833
$ws.Range("I77").Value = 15984
$ws.Range("J77").Value = 10507.25
$ws.Range("K77").Value = 47952
$ws.Range("L77").Value = 31521.75
$ws.Range("M77").Value = -43272
$ws.Range("N77").Value = -40881.75
$ws.Range("H81").Value = 4427.857
$ws.Range("I81").Value = 3589.818
$ws.Range("J81").Value = 5349.7
$ws.Range("K81").Value = 7179.636
$ws.Range("L81").Value = 10699.4
$ws.Range("M81").Value = -6118.636
$ws.Range("N81").Value = -12821.4
$ws.Range("H84").Value = 4427.857
$ws.Range("I84").Value = 3589.818
$ws.Range("J84").Value = 5349.7
$ws.Range("K84").Value = 35898.18
$ws.Range("L84").Value = 53497
$ws.Range("M84").Value = -30594.18
$ws.Range("N84").Value = -64105
$ws.Range("H132").Value = 2711.7256
$ws.Range("I132").Value = 2827.3513
$ws.Range("K132").Value = 8482.053899999999
$ws.Range("M132").Value = -5952.053899999999
$ws.Range("H136").Value = 25002742
$ws.Range("I136").Value = 38462990
$ws.Range("J136").Value = 5135.857
$ws.Range("K136").Value = 115388970
$ws.Range("L136").Value = 15407.571
$ws.Range("M136").Value = -115386420
$ws.Range("N136").Value = -20507.571
